$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new DAMSLTag (col I) / DialogAct (col J)
$updates = @{
    8  = @{ I = "sd"; J = "Statement-non-opinion" }
    12 = @{ I = "sd"; J = "Statement-non-opinion" }
    17 = @{ I = "sd"; J = "Statement-non-opinion" }
    23 = @{ I = "sv"; J = "Statement-opinion" }
    24 = @{ I = "sv"; J = "Statement-opinion" }
    28 = @{ I = "sv"; J = "Statement-opinion" }
    44 = @{ I = "sv"; J = "Statement-opinion" }
    47 = @{ I = "sv"; J = "Statement-opinion" }
    48 = @{ I = "aa"; J = "Agree/Accept" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("I$row").Value = $vals.I
    $ws.Range("J$row").Value = $vals.J
}
